{"js": "// Apply the pt_base v1.1.0 text revisions to the summary table.\nconst body = context.document.body;\n\nconst replacements = [\n  [\"N 2\", \"N 0\"],\n  [\"Summary 2\", \"Summary 0\"],\n  [\"44.9 (10.1) \", \"45 (10) \"],\n  [\"44.6 (10.1) \", \"45 (10) \"],\n  [\"44.8 (10.1) \", \"45 (10) \"],\n  [\"461 (91.1%)\", \"461 (91.11%)\"],\n  [\"444 (89.9%)\", \"444 (89.88%)\"],\n  [\"905 (90.5%)\", \"905 (90.50%)\"]\n];\n\nfor (const [searchText, newText] of replacements) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the pt_base v1.1.0 text revisions to the summary table.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"N 2\"; Replace = \"N 0\" },\n    @{ Find = \"Summary 2\"; Replace = \"Summary 0\" },\n    @{ Find = \"44.9 (10.1) \"; Replace = \"45 (10) \" },\n    @{ Find = \"44.6 (10.1) \"; Replace = \"45 (10) \" },\n    @{ Find = \"44.8 (10.1) \"; Replace = \"45 (10) \" },\n    @{ Find = \"461 (91.1%)\"; Replace = \"461 (91.11%)\" },\n    @{ Find = \"444 (89.9%)\"; Replace = \"444 (89.88%)\" },\n    @{ Find = \"905 (90.5%)\"; Replace = \"905 (90.50%)\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
